# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) counts computed for rows 2..37, replacing the old
# Strike# derived values in column G.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 3
    6  = 3
    7  = 5
    8  = 7
    9  = 5
    10 = 2
    11 = 4
    12 = 2
    13 = 2
    14 = 1
    15 = 6
    16 = 5
    17 = 4
    18 = 2
    19 = 2
    20 = 4
    21 = 6
    22 = 4
    23 = 4
    24 = 3
    25 = 13
    26 = 2
    27 = 7
    28 = 2
    29 = 3
    30 = 3
    31 = 5
    32 = 6
    33 = 3
    34 = 1
    35 = 1
    36 = 4
    37 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
